# MasterQuest.xlsx edit: split the old "_20000_Main001_CallFather" quest row
# into a new "_None" placeholder entry (id 20000) and a renumbered
# "_20001_Main001_CallFather" entry (id 20001), then make the Entities sheet
# the active/visible tab (it was Notes before).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Entities")

# Insert a brand-new row above the old "20000" row; this pushes the existing
# 20000 (Main001 CallFather) and 25000 (Sub001 CallCoin) rows down by one.
$ws1.Rows.Item(2).Insert()

# The row that used to be "20000 / _20000_Main001_CallFather" is now row 3 -
# renumber it to 20001 / _20001_Main001_CallFather (everything else about
# that quest stays the same).
$ws1.Cells.Item(3, 1).Value = 20001
$ws1.Cells.Item(3, 2).Value = "_20001_Main001_CallFather"

# Fill the freshly inserted row 2 with a minimal "_None" placeholder entry.
$ws1.Cells.Item(2, 1).Value = 20000
$ws1.Cells.Item(2, 2).Value = "_None"

# Entities becomes the active sheet (was Notes), with D7 selected.
$ws1.Activate()
[void]$ws1.Range("D7").Select()
